$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting CONTACT NO / Is Mobile? right.
$ws.Columns("C").Insert()

# New column takes the same display width as its left neighbour (LAST NAME).
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# New COUNTRY column header and values.
$ws.Range("C1").Value = "COUNTRY"
$ws.Range("C2").Value = "NL"
$ws.Range("C3").Value = "NL"

# The "Is Mobile?" boolean values (now column E) are no longer populated.
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()

# Update the saved selection to match the authored state.
$ws.Range("E6").Select() | Out-Null
